# Scheduled runner update: refresh computed profit figures on the
# per-job leve/crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with latest market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4343.7407
$ws.Range("I64").Value = 4366.1875
$ws.Range("K64").Value = 4366.1875
$ws.Range("M64").Value = -4118.1875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4343.7407
$ws.Range("I67").Value = 4366.1875
$ws.Range("K67").Value = 4366.1875
$ws.Range("M67").Value = -3508.1875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 10987.079
$ws.Range("I74").Value = 11152.917
$ws.Range("K74").Value = 11152.917
$ws.Range("M74").Value = -10216.917

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 10987.079
$ws.Range("I77").Value = 11152.917
$ws.Range("K77").Value = 55764.585
$ws.Range("M77").Value = -51084.585

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 818.5789
$ws.Range("I96").Value = 584.875
$ws.Range("J96").Value = 988.5454999999999
$ws.Range("K96").Value = 1754.625
$ws.Range("L96").Value = 2965.6365
$ws.Range("M96").Value = -381.625
$ws.Range("N96").Value = -5711.6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2871.2903
$ws.Range("I138").Value = 2285.3044
$ws.Range("J138").Value = 4556
$ws.Range("K138").Value = 6855.9132
$ws.Range("L138").Value = 13668
$ws.Range("M138").Value = -1715.9132
$ws.Range("N138").Value = -23948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1279.2916
$ws.Range("I2").Value = 1672.2307
$ws.Range("J2").Value = 814.9091
$ws.Range("K2").Value = 1672.2307
$ws.Range("L2").Value = 814.9091
$ws.Range("M2").Value = -1559.2307
$ws.Range("N2").Value = -1040.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4027.889
$ws.Range("I32").Value = 2687.5374
$ws.Range("K32").Value = 2687.5374
$ws.Range("M32").Value = -2400.5374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1800.8
$ws.Range("I63").Value = 1800.8
$ws.Range("K63").Value = 1800.8
$ws.Range("M63").Value = -1114.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1800.8
$ws.Range("I66").Value = 1800.8
$ws.Range("K66").Value = 9004
$ws.Range("M66").Value = -5572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1279.2916
$ws.Range("I116").Value = 1672.2307
$ws.Range("J116").Value = 814.9091
$ws.Range("K116").Value = 1672.2307
$ws.Range("L116").Value = 814.9091
$ws.Range("M116").Value = 621.7692999999999
$ws.Range("N116").Value = -5402.9091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1279.2916
$ws.Range("I3").Value = 1672.2307
$ws.Range("J3").Value = 814.9091
$ws.Range("K3").Value = 1672.2307
$ws.Range("L3").Value = 814.9091
$ws.Range("M3").Value = -1558.2307
$ws.Range("N3").Value = -1042.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2677.1177
$ws.Range("I94").Value = 2114.3
$ws.Range("J94").Value = 3481.1428
$ws.Range("K94").Value = 2114.3
$ws.Range("L94").Value = 3481.1428
$ws.Range("M94").Value = -1663.3
$ws.Range("N94").Value = -4383.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 560.64
$ws.Range("I107").Value = 357.17648
$ws.Range("J107").Value = 993
$ws.Range("K107").Value = 357.17648
$ws.Range("L107").Value = 993
$ws.Range("M107").Value = 1562.82352
$ws.Range("N107").Value = -4833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4888.933
$ws.Range("I134").Value = 3766.8635
$ws.Range("J134").Value = 7974.625
$ws.Range("K134").Value = 11300.5905
$ws.Range("L134").Value = 23923.875
$ws.Range("M134").Value = -8765.5905
$ws.Range("N134").Value = -28993.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3931
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3931
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2872.5264
$ws.Range("I68").Value = 6083.3335
$ws.Range("J68").Value = 2270.5
$ws.Range("K68").Value = 18250.0005
$ws.Range("L68").Value = 6811.5
$ws.Range("M68").Value = -17439.0005
$ws.Range("N68").Value = -8433.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2248
$ws.Range("J69").Value = 4014
$ws.Range("L69").Value = 12042
$ws.Range("N69").Value = -13664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2872.5264
$ws.Range("I71").Value = 6083.3335
$ws.Range("J71").Value = 2270.5
$ws.Range("K71").Value = 54750.0015
$ws.Range("L71").Value = 20434.5
$ws.Range("M71").Value = -50694.0015
$ws.Range("N71").Value = -28546.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2248
$ws.Range("J72").Value = 4014
$ws.Range("L72").Value = 36126
$ws.Range("N72").Value = -44238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4129.4
$ws.Range("J80").Value = 4036.75
$ws.Range("L80").Value = 12110.25
$ws.Range("N80").Value = -13982.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4129.4
$ws.Range("J83").Value = 4036.75
$ws.Range("L83").Value = 36330.75
$ws.Range("N83").Value = -45690.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 9176.666999999999
$ws.Range("I139").Value = 2106.75
$ws.Range("K139").Value = 6320.25
$ws.Range("M139").Value = -1180.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5802.0415
$ws.Range("I140").Value = 12463.444
$ws.Range("J140").Value = 1805.2
$ws.Range("K140").Value = 37390.33199999999
$ws.Range("L140").Value = 5415.6
$ws.Range("M140").Value = -32210.33199999999
$ws.Range("N140").Value = -15775.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5417.2
$ws.Range("J80").Value = 4950.9165
$ws.Range("L80").Value = 4950.9165
$ws.Range("N80").Value = -6946.9165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5417.2
$ws.Range("J83").Value = 4950.9165
$ws.Range("L83").Value = 24754.5825
$ws.Range("N83").Value = -34738.5825

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3403.2258
$ws.Range("I132").Value = 2868.4385
$ws.Range("K132").Value = 8605.315500000001
$ws.Range("M132").Value = -6075.315500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2348.1667
$ws.Range("J16").Value = 5444.75
$ws.Range("L16").Value = 5444.75
$ws.Range("N16").Value = -5784.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1610.7646
$ws.Range("I22").Value = 1039.5385
$ws.Range("J22").Value = 3467.25
$ws.Range("K22").Value = 1039.5385
$ws.Range("L22").Value = 3467.25
$ws.Range("M22").Value = -744.5385000000001
$ws.Range("N22").Value = -4057.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1610.7646
$ws.Range("I27").Value = 1039.5385
$ws.Range("J27").Value = 3467.25
$ws.Range("K27").Value = 1039.5385
$ws.Range("L27").Value = 3467.25
$ws.Range("M27").Value = -932.5385000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3410.8235
$ws.Range("I68").Value = 3414.6
$ws.Range("J68").Value = 3407.842
$ws.Range("K68").Value = 3414.6
$ws.Range("L68").Value = 3407.842
$ws.Range("M68").Value = -2665.6
$ws.Range("N68").Value = -4905.842000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3410.8235
$ws.Range("I71").Value = 3414.6
$ws.Range("J71").Value = 3407.842
$ws.Range("K71").Value = 17073
$ws.Range("L71").Value = 17039.21
$ws.Range("M71").Value = -13329
$ws.Range("N71").Value = -24527.21

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 708.0625
$ws.Range("I113").Value = 427.16666
$ws.Range("K113").Value = 1281.49998
$ws.Range("M113").Value = 888.5000199999999
